# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45719 (2025-03-03) to 45720 (2025-03-04).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 42; $row++) {
    $ws.Cells.Item($row, 3).Value = 45720
}
